# Corrected excel sheets for application fix issues
#
# Updates the Summary, Repayment schedule and Transactions sheets to
# reflect a corrected repayment schedule (an extra "Over Due" column
# on the schedule, updated due/interest figures, and renumbered
# transaction ids), and leaves the Transactions sheet as the active
# sheet/selection when the workbook is saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("F2").Value = 0
$summary.Range("A3").Value = 213.52
$summary.Range("E3").Value = 113.52

# ---------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Repayment schedule")
$schedule.Activate() | Out-Null

# Introduce the new "O" (Over Due) column for every data row, copying
# the formatting of the neighbouring "N" column so the new cells pick
# up the same style as the rest of the sheet.
$schedule.Range("N2").Copy($schedule.Range("O2")) | Out-Null

$schedule.Range("N3").Copy($schedule.Range("O3")) | Out-Null
$schedule.Range("O3").Value = 0

$schedule.Range("N4").Copy($schedule.Range("O4")) | Out-Null
$schedule.Range("O4").Value = 0

$schedule.Range("N5").Copy($schedule.Range("O5")) | Out-Null
$schedule.Range("O5").Value = 0

$schedule.Range("N6").Copy($schedule.Range("O6")) | Out-Null
$schedule.Range("O6").Value = 0

$schedule.Range("N7").Copy($schedule.Range("O7")) | Out-Null
$schedule.Range("O7").Value = 0

$schedule.Range("N8").Copy($schedule.Range("O8")) | Out-Null
$schedule.Range("O8").Value = 0

# Row 4 (installment 2)
$schedule.Range("C4").Value = 42095
$schedule.Range("F4").Value = 923.19
$schedule.Range("G4").Value = 3212.27
$schedule.Range("H4").Value = 41.35

# Row 5 (installment 3)
$schedule.Range("B5").Value = 30
$schedule.Range("C5").Value = 42125
$schedule.Range("F5").Value = 932.42
$schedule.Range("G5").Value = 2279.85
$schedule.Range("H5").Value = 32.12

# Row 6 (installment 4)
$schedule.Range("B6").Value = 31
$schedule.Range("C6").Value = 42156
$schedule.Range("F6").Value = 941.74
$schedule.Range("G6").Value = 1338.11
$schedule.Range("H6").Value = 22.8

# Row 7 (installment 5)
$schedule.Range("B7").Value = 30
$schedule.Range("C7").Value = 42186
$schedule.Range("F7").Value = 951.16
$schedule.Range("G7").Value = 386.95
$schedule.Range("H7").Value = 13.38

# Row 8 (installment 6 - final)
$schedule.Range("B8").Value = 31
$schedule.Range("C8").Value = 42217
$schedule.Range("F8").Value = 386.95
$schedule.Range("H8").Value = 3.87
$schedule.Range("K8").Value = 390.82
$schedule.Range("P8").Value = 390.82

# Schedule sheet is no longer the tab shown when the workbook is
# reopened - select the row below the data instead of the old L8 cell.
$schedule.Range("A9:XFD9").Select() | Out-Null

# ---------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------
$transactions = $wb.Worksheets.Item("Transactions")
$transactions.Range("A2").Value = 83
$transactions.Range("A3").Value = 82

# Transactions becomes the active/selected sheet and range.
$transactions.Activate() | Out-Null
$transactions.Range("A2:L3").Select() | Out-Null
